$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 2945.0527
$ws.Range("I64").Value = 2809.6365
$ws.Range("K64").Value = 2809.6365
$ws.Range("M64").Value = -2561.6365
# Row 67
$ws.Range("H67").Value = 2945.0527
$ws.Range("I67").Value = 2809.6365
$ws.Range("K67").Value = 2809.6365
$ws.Range("M67").Value = -1951.6365
# Row 82
$ws.Range("H82").Value = 600
$ws.Range("I82").Value = 600
$ws.Range("K82").Value = 1800
$ws.Range("M82").Value = -1394
# Row 85
$ws.Range("H85").Value = 600
$ws.Range("I85").Value = 600
$ws.Range("K85").Value = 1800
$ws.Range("M85").Value = -396
# Row 98
$ws.Range("H98").Value = 474.72223
$ws.Range("I98").Value = 462.08334
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 462.08334
$ws.Range("L98").Value = 500
$ws.Range("M98").Value = 1035.91666
$ws.Range("N98").Value = -3496
# Row 122
$ws.Range("H122").Value = 474.72223
$ws.Range("I122").Value = 462.08334
$ws.Range("J122").Value = 500
$ws.Range("K122").Value = 1386.25002
$ws.Range("L122").Value = 1500
$ws.Range("M122").Value = 1063.74998
$ws.Range("N122").Value = -6400
# Row 125
$ws.Range("H125").Value = 1723.1428
$ws.Range("I125").Value = 1612.4
$ws.Range("K125").Value = 14511.6
$ws.Range("M125").Value = -12051.6
# Row 137
$ws.Range("H137").Value = 1496.1476
$ws.Range("I137").Value = 1329.8379
$ws.Range("J137").Value = 1752.5416
$ws.Range("K137").Value = 3989.5137
$ws.Range("L137").Value = 5257.6248
$ws.Range("M137").Value = -1439.5137
$ws.Range("N137").Value = -10357.6248

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8020.846
$ws.Range("I32").Value = 2575.3333
$ws.Range("J32").Value = 25099.955
$ws.Range("K32").Value = 2575.3333
$ws.Range("L32").Value = 25099.955
$ws.Range("M32").Value = -2288.3333
$ws.Range("N32").Value = -25673.955
# Row 61
$ws.Range("H61").Value = 2327.6072
$ws.Range("I61").Value = 1994.6923
$ws.Range("J61").Value = 2616.1333
$ws.Range("K61").Value = 1994.6923
$ws.Range("L61").Value = 2616.1333
$ws.Range("M61").Value = -1782.6923
$ws.Range("N61").Value = -3040.1333
# Row 136
$ws.Range("H136").Value = 2327.6072
$ws.Range("I136").Value = 1994.6923
$ws.Range("J136").Value = 2616.1333
$ws.Range("K136").Value = 5984.0769
$ws.Range("L136").Value = 7848.3999
$ws.Range("M136").Value = -3434.0769
$ws.Range("N136").Value = -12948.3999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 68
$ws.Range("H68").Value = 43000
$ws.Range("J68").Value = 43000
$ws.Range("L68").Value = 43000
$ws.Range("N68").Value = -44622
# Row 71
$ws.Range("H71").Value = 43000
$ws.Range("J71").Value = 43000
$ws.Range("L71").Value = 129000
$ws.Range("N71").Value = -137112

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 42000
$ws.Range("J20").Value = 42000
$ws.Range("L20").Value = 42000
$ws.Range("N20").Value = -42472
# Row 30
$ws.Range("H30").Value = 42000
$ws.Range("J30").Value = 42000
$ws.Range("L30").Value = 42000
$ws.Range("N30").Value = -42182
# Row 31
$ws.Range("H31").Value = 5557266
$ws.Range("I31").Value = 10527044
$ws.Range("J31").Value = 3775647.8
$ws.Range("K31").Value = 10527044
$ws.Range("L31").Value = 3775647.8
$ws.Range("M31").Value = -10526749
$ws.Range("N31").Value = -3776237.8
# Row 34
$ws.Range("H34").Value = 5557266
$ws.Range("I34").Value = 10527044
$ws.Range("J34").Value = 3775647.8
$ws.Range("K34").Value = 10527044
$ws.Range("L34").Value = 3775647.8
$ws.Range("M34").Value = -10526842
$ws.Range("N34").Value = -3776051.8
# Row 58
$ws.Range("H58").Value = 2779
$ws.Range("I58").Value = 933.7895
$ws.Range("J58").Value = 5700.5835
$ws.Range("K58").Value = 933.7895
$ws.Range("L58").Value = 5700.5835
$ws.Range("M58").Value = -730.7895
$ws.Range("N58").Value = -6106.5835
# Row 62
$ws.Range("H62").Value = 166669550
$ws.Range("J62").Value = 200002740
$ws.Range("L62").Value = 200002740
$ws.Range("N62").Value = -200003988
# Row 65
$ws.Range("H65").Value = 166669550
$ws.Range("J65").Value = 200002740
$ws.Range("L65").Value = 1000013700
$ws.Range("N65").Value = -1000019940
# Row 128
$ws.Range("H128").Value = 42000
$ws.Range("J128").Value = 42000
$ws.Range("L128").Value = 42000
$ws.Range("N128").Value = -51960
# Row 129
$ws.Range("H129").Value = 26785.428
$ws.Range("J129").Value = 26785.428
$ws.Range("L129").Value = 26785.428
$ws.Range("N129").Value = -36785.428
# Row 136
$ws.Range("H136").Value = 2779
$ws.Range("I136").Value = 933.7895
$ws.Range("J136").Value = 5700.5835
$ws.Range("K136").Value = 2801.3685
$ws.Range("L136").Value = 17101.7505
$ws.Range("M136").Value = -251.3685
$ws.Range("N136").Value = -22201.7505

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 747.07367
$ws.Range("I68").Value = 425.8154
$ws.Range("J68").Value = 1443.1333
$ws.Range("K68").Value = 1277.4462
$ws.Range("L68").Value = 4329.3999
$ws.Range("M68").Value = -466.4462000000001
$ws.Range("N68").Value = -5951.3999
# Row 71
$ws.Range("H71").Value = 747.07367
$ws.Range("I71").Value = 425.8154
$ws.Range("J71").Value = 1443.1333
$ws.Range("K71").Value = 3832.3386
$ws.Range("L71").Value = 12988.1997
$ws.Range("M71").Value = 223.6614
$ws.Range("N71").Value = -21100.1997
# Row 92
$ws.Range("H92").Value = 199.88889
$ws.Range("I92").Value = 162.375
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 487.125
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 760.875
$ws.Range("N92").Value = -3996
# Row 98
$ws.Range("H98").Value = 597
$ws.Range("J98").Value = 597
$ws.Range("L98").Value = 1791
$ws.Range("N98").Value = -4787
# Row 129
$ws.Range("H129").Value = 7938189.5
$ws.Range("I129").Value = 1147.6
$ws.Range("J129").Value = 27780794
$ws.Range("K129").Value = 3442.8
$ws.Range("L129").Value = 83342382
$ws.Range("M129").Value = 1557.2
$ws.Range("N129").Value = -83352382

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2391.6785
$ws.Range("I132").Value = 1653.8334
$ws.Range("K132").Value = 4961.5002
$ws.Range("M132").Value = -2431.5002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 19904176
$ws.Range("I68").Value = 28196110
$ws.Range("J68").Value = 3530.6
$ws.Range("K68").Value = 28196110
$ws.Range("L68").Value = 3530.6
$ws.Range("M68").Value = -28195361
$ws.Range("N68").Value = -5028.6
# Row 71
$ws.Range("H71").Value = 19904176
$ws.Range("I71").Value = 28196110
$ws.Range("J71").Value = 3530.6
$ws.Range("K71").Value = 140980550
$ws.Range("L71").Value = 17653
$ws.Range("M71").Value = -140976806
$ws.Range("N71").Value = -25141
# Row 122
$ws.Range("H122").Value = 4669.9287
$ws.Range("J122").Value = 5108.3335
$ws.Range("L122").Value = 15325.0005
$ws.Range("N122").Value = -20225.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2535.5833
$ws.Range("I132").Value = 1992.1111
$ws.Range("J132").Value = 4166
$ws.Range("K132").Value = 5976.3333
$ws.Range("L132").Value = 12498
$ws.Range("M132").Value = -3446.3333
$ws.Range("N132").Value = -17558
